$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One-off data correction: NOSH for row 364 (DE... ticker) was mis-keyed.
$ws.Range("G364").Value = 700651000

# Remove the "WALLENIUS WILHELMS" (WAWI-NO / NO0010571680) entry from the
# Universe table entirely. Deleting the whole row shifts every row below it
# up by one and shrinks the table/used range from A1:M438 to A1:M437.
$ws.Rows(427).Delete()
